# New PO forecast model
# - Weekly Quantity: append one more week (row 5)
# - Monthly Trend: append one more month (row 5)
# - PO Forecast: re-forecast the weekly series starting 2025-01-26,
#   replacing rows 5-12 and adding a new row 13, each with qty 1.

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---- Sheet 1: Weekly Quantity ----
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A5").Value = 45683.99999999999
$wsWeekly.Range("A5").NumberFormat = $dateFmt
$wsWeekly.Range("B5").Value = 1

# ---- Sheet 2: Monthly Trend ----
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A5").Value = 45688.99999999999
$wsMonthly.Range("A5").NumberFormat = $dateFmt
$wsMonthly.Range("B5").Value = 1

# ---- Sheet 3: PO Forecast ----
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastDates = @(45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999, 45732.99999999999, 45739.99999999999)

$row = 5
foreach ($d in $forecastDates) {
    $aCell = $wsForecast.Cells.Item($row, 1)
    $aCell.Value = $d
    $aCell.NumberFormat = $dateFmt
    $bCell = $wsForecast.Cells.Item($row, 2)
    $bCell.Value = 1
    $row = $row + 1
}
